$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert 9 new rows (26-34) for additional chb01 Test/Interictal files,
# right after the existing chb01_08.edf row (row 25).
$ws.Range("A26:A34").EntireRow.Insert()

# Insert 5 new rows (41-45) for additional chb06 Test/Interictal files,
# right after the existing chb06_08.edf row (now row 40).
$ws.Range("A41:A45").EntireRow.Insert()

# Fill the new chb01 rows.
$chb01Files = @("chb01_09.edf","chb01_10.edf","chb01_11.edf","chb01_12.edf","chb01_13.edf","chb01_14.edf","chb01_17.edf","chb01_19.edf","chb01_20.edf")
for ($i = 0; $i -lt $chb01Files.Length; $i++) {
    $r = 26 + $i
    $ws.Cells.Item($r,1).Value = "chb01"
    $ws.Cells.Item($r,2).Value = $chb01Files[$i]
    $ws.Cells.Item($r,3).Value = "Test"
    $ws.Cells.Item($r,4).Value = "Interictal"
    $ws.Cells.Item($r,5).Value = 0
    $ws.Cells.Item($r,6).Value = 0
    $ws.Cells.Item($r,7).Value = 900
}

# Fill the new chb06 rows. The last two (chb06_16.edf, chb06_17.edf) get
# double the crop length, expressed as a formula (900*2), matching the
# other extended-length rows already in the sheet.
$chb06Files = @("chb06_12.edf","chb06_14.edf","chb06_15.edf","chb06_16.edf","chb06_17.edf")
for ($i = 0; $i -lt $chb06Files.Length; $i++) {
    $r = 41 + $i
    $ws.Cells.Item($r,1).Value = "chb06"
    $ws.Cells.Item($r,2).Value = $chb06Files[$i]
    $ws.Cells.Item($r,3).Value = "Test"
    $ws.Cells.Item($r,4).Value = "Interictal"
    $ws.Cells.Item($r,5).Value = 0
    $ws.Cells.Item($r,6).Value = 0
    if ($r -ge 44) {
        $ws.Cells.Item($r,7).Formula = "=900*2"
    } else {
        $ws.Cells.Item($r,7).Value = 900
    }
}

# Re-apply the worksheet AutoFilter over the now-larger table and keep the
# underlying defined-name/filter-database hidden, like Excel does natively.
$ws.Range("A1:G51").AutoFilter()
$flt = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$G`$51")
$flt.Visible = $false

# Scroll / select like the saved file did after the edit.
$ws.Application.Goto($ws.Range("A14"))
$ws.Range("A45").Select()

# Restore window size/position (cosmetic, matches the saved workbook view).
$excel.ActiveWindow.WindowState = -4143
$excel.Width = 21920
$excel.Height = 19920
$excel.Left = 540
$excel.Top = 500

Write-Host "done"
